$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.555.17"
$ws.Range("E2").Value = "  +0.18%  "
$ws.Range("D3").Value = "3.141.94"
$ws.Range("E3").Value = "  -1.18%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.07"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "164.69"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.76%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.576"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.48%  "
$ws.Range("D9").Value = "3.159.46"
$ws.Range("E9").Value = "  -0.94%  "
$ws.Range("E10").Value = "  -2.44%  "
$ws.Range("E11").Value = "  -2.58%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.385"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.20%  "
$ws.Range("D13").Value = "3.690.50"
$ws.Range("E13").Value = "  -1.40%  "
$ws.Range("E14").Value = "  -1.66%  "
$ws.Range("D15").Value = "64.593.87"
$ws.Range("E15").Value = "  +0.13%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.01"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.21%  "
$ws.Range("D17").Value = "3.149.12"
$ws.Range("E17").Value = "  -1.33%  "
$ws.Range("E18").Value = "  -2.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "408.86"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.78%  "
$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.60%  "
$ws.Range("B21").Value = "Polkadot"
$ws.Range("C21").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.08"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.31%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.88"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.00%  "
$ws.Range("E25").Value = "  -2.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.194"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.15%  "
$ws.Range("E27").Value = "  -2.63%  "
$ws.Range("E28").Value = "  +1.83%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.997"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.07%  "
$ws.Range("E31").Value = "  -1.05%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.26"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.54%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "163.35"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.00%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.89"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.99%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.30"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.89%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.13"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.07%  "
$ws.Range("E37").Value = "  -0.29%  "
$ws.Range("E38").Value = "  -0.97%  "
$ws.Range("D39").Value = "2.636.69"
$ws.Range("E39").Value = "  -2.67%  "
$ws.Range("E40").Value = "  -2.24%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.10"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.89%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.22"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.76%  "
$ws.Range("E43").Value = "  -3.33%  "
$ws.Range("E44").Value = "  -1.26%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.33"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.11%  "
$ws.Range("B46").Value = "Bittensor"
$ws.Range("C46").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "291.34"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.18%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "21.40"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("E48").Value = "  -3.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.997"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0975"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.24%  "
$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "10.48"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.44%  "
